# Refresh the cryptocurrency price / volume table with the latest scrape results.
# (Mirrors the "Updated cryptos list ... with GitHub Actions" automated commit.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores numeric-looking price text (e.g. "1.00", "542.81") as plain text
# in the workbook (inline strings, not numbers). Setting .Value on a numeric-looking
# string would make Excel auto-convert it to a real number, so pre-format just the
# cells that need this protection as Text ("@") before writing their new value.
$textPriceCells = @("D5", "D6", "D7", "D10", "D12", "D13", "D18", "D19", "D20", "D21", "D24", "D25", "D26", "D27", "D30", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textPriceCells) { $ws.Range($addr).NumberFormat = "@" }

# Write every updated cell (Coin name / Link / Price / Volume(1h)).
$ws.Range("D2").Value = "58.359.46"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "2.349.14"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "542.81"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "135.05"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E8").Value = "  +5.20%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "5.67"
$ws.Range("E10").Value = "  +6.58%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("D13").Value = "23.79"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "2.767.63"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "58.288.90"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "2.327.44"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "10.73"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").Value = "333.29"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").Value = "4.26"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").Value = "  -3.26%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("D24").Value = "62.77"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "8.50"
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  +5.48%  "
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("D30").Value = "170.35"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "0.0₃0736"
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "18.42"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("B34").Value = "SuiNetwork"
$ws.Range("C34").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D34").Value = "1.02"
$ws.Range("E34").Value = "  +12.11%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "4.25"
$ws.Range("E36").Value = "  +5.75%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").Value = "1.64"
$ws.Range("E39").Value = "  +3.59%  "
$ws.Range("D40").Value = "39.16"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "143.27"
$ws.Range("E41").Value = "  -3.43%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "0.377"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.64"
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "288.43"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").Value = "0.0942"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "19.19"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("D47").Value = "0.0503"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "0.564"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "0.382"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "17.51"
$ws.Range("E51").Value = "  +0.79%  "
